$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.003.84"
$ws.Range("E2").Value = "  -3.28%  "

$ws.Range("D3").Value = "1.714.13"
$ws.Range("E3").Value = "  -3.03%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'308.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.02%  "

$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("E7").Value = "  +4.68%  "

$ws.Range("D8").Value = "'0.3451"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.11%  "

$ws.Range("D9").Value = "'41.88"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.36%  "

$ws.Range("D10").Value = "'0.07236"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.99%  "

$ws.Range("E11").Value = "  -5.26%  "

$ws.Range("E12").Value = "  +0.10%  "

$ws.Range("E13").Value = "  -4.66%  "

$ws.Range("D14").Value = "'5.814"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.25%  "

$ws.Range("D15").Value = "1.718.47"
$ws.Range("E15").Value = "  -3.29%  "

$ws.Range("D16").Value = "'6.816"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.10%  "

$ws.Range("D17").Value = "'86.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.35%  "

$ws.Range("E18").Value = "  -2.53%  "

$ws.Range("D19").Value = "'0.06379"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.08%  "

$ws.Range("D20").Value = "'1.002"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").Value = "'16.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.94%  "

$ws.Range("D22").Value = "'5.599"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.87%  "

$ws.Range("D23").Value = "27.060.97"
$ws.Range("E23").Value = "  -3.16%  "

$ws.Range("E24").Value = "  -4.15%  "

$ws.Range("D25").Value = "'2.095"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.19%  "

$ws.Range("D26").Value = "'19.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.14%  "

$ws.Range("D27").Value = "'150.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.52%  "

$ws.Range("D28").Value = "1.913.91"
$ws.Range("E28").Value = "  -3.20%  "

$ws.Range("D29").Value = "'2.061"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.15%  "

$ws.Range("D30").Value = "'120.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.02%  "

$ws.Range("D31").Value = "'1.024"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.75%  "

$ws.Range("D32").Value = "'0.09166"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.15%  "

$ws.Range("D33").Value = "'3.600"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.77%  "

$ws.Range("D34").Value = "'5.294"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.67%  "

$ws.Range("D35").Value = "'1.472"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.44%  "

$ws.Range("D36").Value = "'0.02169"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.06%  "

$ws.Range("D37").Value = "'0.05824"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.75%  "

$ws.Range("D38").Value = "'0.1992"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.79%  "

$ws.Range("D39").Value = "'10.87"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.07%  "

$ws.Range("D41").Value = "'4.694"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.00%  "

$ws.Range("D42").Value = "'0.5953"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.76%  "

$ws.Range("E43").Value = "  -7.99%  "

$ws.Range("D44").Value = "'7.476"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.08%  "

$ws.Range("D45").Value = "'12.69"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.96%  "

$ws.Range("D46").Value = "'3.573"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.39%  "

$ws.Range("D47").Value = "'0.5558"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.85%  "

$ws.Range("D48").Value = "'118.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.97%  "

$ws.Range("D49").Value = "'1.818"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.85%  "

$ws.Range("E50").Value = "  -2.11%  "

$ws.Range("D51").Value = "'0.06627"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.95%  "
